$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the existing date-format style used in column D (e.g. cell D136)
$dateFormat = $ws.Cells.Item(136, 4).NumberFormat

# Row 137
$ws.Cells.Item(137, 1).Value = 3
$ws.Cells.Item(137, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(137, 3).Value = "Coquimbo"
$ws.Cells.Item(137, 4).NumberFormat = $dateFormat
$ws.Cells.Item(137, 4).Value = 44890
$ws.Cells.Item(137, 5).Value = 5
$ws.Cells.Item(137, 6).Value = "Fruta"
$ws.Cells.Item(137, 7).Value = 100103
$ws.Cells.Item(137, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(137, 9).Value = 100103003
$ws.Cells.Item(137, 10).Value = "Damasco"
$ws.Cells.Item(137, 11).Value = "Castle Brite"
$ws.Cells.Item(137, 12).Value = "Especial"
$ws.Cells.Item(137, 13).Value = 50
$ws.Cells.Item(137, 14).Value = 15000
$ws.Cells.Item(137, 15).Value = 15000
$ws.Cells.Item(137, 16).Value = 15000
$ws.Cells.Item(137, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(137, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(137, 19).Value = 1500
$ws.Cells.Item(137, 20).Value = 10

# Row 138
$ws.Cells.Item(138, 1).Value = 3
$ws.Cells.Item(138, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(138, 3).Value = "Coquimbo"
$ws.Cells.Item(138, 4).NumberFormat = $dateFormat
$ws.Cells.Item(138, 4).Value = 44890
$ws.Cells.Item(138, 5).Value = 5
$ws.Cells.Item(138, 6).Value = "Fruta"
$ws.Cells.Item(138, 7).Value = 100103
$ws.Cells.Item(138, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(138, 9).Value = 100103003
$ws.Cells.Item(138, 10).Value = "Damasco"
$ws.Cells.Item(138, 11).Value = "Castle Brite"
$ws.Cells.Item(138, 12).Value = "Primera"
$ws.Cells.Item(138, 13).Value = 57
$ws.Cells.Item(138, 14).Value = 13000
$ws.Cells.Item(138, 15).Value = 13000
$ws.Cells.Item(138, 16).Value = 13000
$ws.Cells.Item(138, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(138, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(138, 19).Value = 1300
$ws.Cells.Item(138, 20).Value = 10
